# Rename the paired "_old"/"_new" column headers to the format-version based
# names "_FV2410" (old / "before") and "_FV2504" (new / "after"), add a
# worksheet Table (ListObject) over the used range with an AutoFilter, and
# freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row (row 1) cells --------------------------------
$headers = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410",
    "diff",
    "Segmentname_FV2504",
    "Segmentgruppe_FV2504",
    "Segment_FV2504",
    "Datenelement_FV2504",
    "Segment ID_FV2504",
    "Code_FV2504",
    "Qualifier_FV2504",
    "Beschreibung_FV2504",
    "Bedingungsausdruck_FV2504",
    "Bedingung_FV2504"
)

for ($i = 0; $i -lt $headers.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value2 = $headers[$i]
}

# --- 2. Turn the used range into a proper Excel Table (ListObject) -----
$usedRange = $ws.Range("A1:U77")
$lo = $ws.ListObjects.Add(1, $usedRange, [System.Reflection.Missing]::Value, 1)
$lo.Name = "Table1"

# --- 3. Freeze the header row -------------------------------------------
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select() | Out-Null
